# Apply the "Issue - #9:Done, #17:Done, #18 추가" edit to the workbook.
$wb = $excel.ActiveWorkbook

$wsIssue = $wb.Worksheets.Item("Issue list")
$wsEco   = $wb.Worksheets.Item("ECO list")

# ---------------------------------------------------------------------------
# Workbook-level: recalculation becomes automatic again, and the
# "Issue list" sheet becomes the active / selected tab.
# ---------------------------------------------------------------------------
$wb.Application.Calculation = -4105   # xlCalculationAutomatic
$wsIssue.Activate()

# ---------------------------------------------------------------------------
# Sheet "Issue list"
# ---------------------------------------------------------------------------

# --- Row 12 (Issue #9) : now finished / closed -----------------------------
# Pull formatting from a sibling "normal" data row (row 19) which uses the
# exact same style set as the target row, but keep row 12 at its existing
# (default) row height.
$wsIssue.Range("B19:I19").Copy()
$wsIssue.Range("B12:I12").PasteSpecial(-4122)   # xlPasteFormats
$wsIssue.Application.CutCopyMode = $false

$wsIssue.Range("B12").Value = 9
$wsIssue.Range("C12").Value = [datetime]"2017-12-30"
$wsIssue.Range("D12").Value = "SW"
$wsIssue.Range("E12").Value = "CLOSED"
$wsIssue.Range("F12").Value = [datetime]"2018-01-03"
$wsIssue.Range("G12").Value = "전원 OFF시 buzzer 동작하지 않음 - 기능 미구현 상태"
$wsIssue.Range("H12").Value = "기능 구현"
$wsIssue.Range("I12").Value = ""

# --- Row 20 (Issue #17) : now finished / closed -----------------------------
$wsIssue.Range("B19:I19").Copy()
$wsIssue.Range("B20:I20").PasteSpecial(-4122)   # xlPasteFormats
$wsIssue.Application.CutCopyMode = $false

$wsIssue.Range("B20").Value = 17
$wsIssue.Range("C20").Value = [datetime]"2017-12-31"
$wsIssue.Range("D20").Value = "SW"
$wsIssue.Range("E20").Value = "CLOSED"
$wsIssue.Range("F20").Value = [datetime]"2018-01-06"
$wsIssue.Range("G20").Value = "Low Battery시 LED_R Blink 동작하지 않음"
$wsIssue.Range("H20").Value = "Chargetask_lowbatt.c / Chargetask_state.h 수정 "
$wsIssue.Range("I20").Value = ""

# --- Row 21 (Issue #18) : brand-new open issue ------------------------------
# Formatting for row 21 already matches what is needed (it is still the
# blank-template row style), only the row height needs to grow to fit the
# wrapped, multi-line description.
$wsIssue.Range("C21").Value = [datetime]"2018-01-05"
$wsIssue.Range("E21").Value = "OPEN"
$newLine = [char]10
$wsIssue.Range("G21").Value = "RS-232 동작 안됨" + $newLine + " - RX interrupt 인식을 못함" + $newLine + " - USART1_RX date의 Low level이 2.17V 까지밖에 안떨어짐(정상 0V)" + $newLine + " - MAX3232의 ChargePump의 파형이 점검치구 board와 틀림"
$wsIssue.Rows(21).RowHeight = 66

# --- sheetView : freeze-pane top-left cell and active selection -----------
$wsIssue.Activate()
$wsIssue.Range("A7").Select()
$excel.ActiveWindow.FreezePanes = $true
$wsIssue.Range("G26").Select()

# ---------------------------------------------------------------------------
# Sheet "ECO list" : active selection moves, and it is no longer the
# selected tab (handled above by activating "Issue list" last… but we must
# still set the ECO sheet's own remembered selection).
# ---------------------------------------------------------------------------
$wsEco.Range("D20").Select()

# Make sure "Issue list" ends up the active sheet/tab (must be last).
$wsIssue.Activate()
$wsIssue.Range("G26").Select()

Write-Host "edit applied"
